# Commit: started changeover from blocking polled SD card writes to data
# portion interrupt driven and non-blocking polling code.
#
# Pinout Translation sheet: the SD card now drives its own power-enable
# line instead of using the PGED/PREC (reversed) debug pins, so clear the
# old "PGED (reversed)" note and relabel the RA1 connection as the new
# "SD Power Enable" line. Also move the active selection to reflect where
# the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# RA0 (row 15) no longer has a "PGED (reversed)" note.
$ws.Range("C15").ClearContents()

# RA1 (row 16) is repurposed from "PREC (reversed)" to the new SD card
# power-enable signal.
$ws.Range("C16").Value = "SD Power Enable"

# Reflect the author's active selection at the time of the edit.
$ws.Range("G19").Select()
